# This script applies the commit's change: the data rows 2 and 3 in the
# worksheet were swapped (row 2 now holds what used to be row 3's record,
# and vice versa), while row numbers / surrounding structure stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Cell, $Text)
    # Force text interpretation so Excel does not auto-convert values that
    # look like numbers/dates (e.g. "1988-01-01", "1") into numeric/date
    # serials - the source data keeps these as plain strings.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
}

function Clear-Cell {
    param($Cell)
    $Cell.ClearContents()
}

# ---- Row 2 gets the values that used to live in row 3 ----
$ws.Cells.Item(2,1).Value = 74601266          # A2
$ws.Cells.Item(2,2).Value = 95519             # B2
$ws.Cells.Item(2,5).Value = 221945            # E2
$ws.Cells.Item(2,6).Value = "Revlummer"       # F2
$ws.Cells.Item(2,7).Value = "Lycopodium annotinum" # G2
$ws.Cells.Item(2,8).Value = "L."              # H2
Clear-Cell $ws.Cells.Item(2,9)                # I2 (now blank)
Clear-Cell $ws.Cells.Item(2,12)               # L2 (no longer present)
Set-TextCell $ws.Cells.Item(2,16) "Pettersborg 200 m NNV t N, Sm" # P2
$ws.Cells.Item(2,17).Value = 579547.1134967525   # Q2
$ws.Cells.Item(2,18).Value = 6427974.619960002   # R2
$ws.Cells.Item(2,19).Value = 50               # S2
Set-TextCell $ws.Cells.Item(2,25) "1988-01-01" # Y2
Set-TextCell $ws.Cells.Item(2,27) "1989-12-31" # AA2
Set-TextCell $ws.Cells.Item(2,29) "Smålands flora 2007: KOO: 7G5g 4219. SOM: Lycopodium annotinum. LEG: Birger Danielsson" # AC2
Set-TextCell $ws.Cells.Item(2,35) "Skogsmossekant" # AI2 (newly present)
Set-TextCell $ws.Cells.Item(2,49) "Margareta Edqvist" # AW2
Set-TextCell $ws.Cells.Item(2,50) "Via Margareta Edqvist" # AX2
Set-TextCell $ws.Cells.Item(2,51) "Smålands flora (1978-2007)" # AY2

# ---- Row 3 gets the values that used to live in row 2 ----
$ws.Cells.Item(3,1).Value = 16879685          # A3
$ws.Cells.Item(3,2).Value = 9302              # B3
$ws.Cells.Item(3,5).Value = 101246            # E3
$ws.Cells.Item(3,6).Value = "Ekoxe"           # F3
$ws.Cells.Item(3,7).Value = "Lucanus cervus"  # G3
$ws.Cells.Item(3,8).Value = "(Linnaeus, 1758)" # H3
Set-TextCell $ws.Cells.Item(3,9) "1"          # I3 (now "1")
Set-TextCell $ws.Cells.Item(3,12) "hane"      # L3 (newly present)
Set-TextCell $ws.Cells.Item(3,16) "Överum, Sm" # P3
$ws.Cells.Item(3,17).Value = 579616.5892805457   # Q3
$ws.Cells.Item(3,18).Value = 6427821.143954458   # R3
$ws.Cells.Item(3,19).Value = 25               # S3
Set-TextCell $ws.Cells.Item(3,25) "2013-07-13" # Y3
Set-TextCell $ws.Cells.Item(3,27) "2013-07-13" # AA3
Set-TextCell $ws.Cells.Item(3,29) "Thore Belinder. 073-3103409. Lst Kan kontakta honom om intresse finns. Finns på plats till 25-07. Ev längre. Medobservatör: Christine Strand" # AC3
Clear-Cell $ws.Cells.Item(3,35)               # AI3 (no longer present)
Set-TextCell $ws.Cells.Item(3,49) "Stefan Karlsson" # AW3
Set-TextCell $ws.Cells.Item(3,50) "Via Stefan Karlsson" # AX3
Set-TextCell $ws.Cells.Item(3,51) "Ekoxeuppropet 2013" # AY3
